# Generate Report for Archive
#
# The localization status for the 5825bb2e... report entries moved on from
# "Ready for handoff" to "In Translation". That phrase is a shared string
# used by the Status column on every per-language sheet (zh-cn, de-de) and
# mirrored onto the Overview roll-up sheet (zh-cn / de-de columns), so
# updating the three status cells on each sheet refreshes all of them.
# The Status/roll-up columns are also narrowed to re-fit the new, shorter
# text.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value2 = $newStatus
$wsOverview.Range("F2").Value2 = $newStatus
$wsOverview.Range("E3").Value2 = $newStatus
$wsOverview.Range("F3").Value2 = $newStatus
$wsOverview.Range("E4").Value2 = $newStatus
$wsOverview.Range("F4").Value2 = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value2 = $newStatus
$wsZhCn.Range("C3").Value2 = $newStatus
$wsZhCn.Range("C4").Value2 = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value2 = $newStatus
$wsDeDe.Range("C3").Value2 = $newStatus
$wsDeDe.Range("C4").Value2 = $newStatus

# Re-fit the columns that held the status text so they hug the shorter
# "In Translation" string instead of the wider "Ready for handoff".
$newWidth = 12.576851254417766
$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
